$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added to the dataset. It belongs right above
# the row that currently holds the oldest record (row 283), so insert a
# fresh row there; Excel shifts rows 283:369 down to 284:370 and the sheet's
# dimension grows from A1:T369 to A1:T370.
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new observation.
$ws.Range('A283').Value = 1
$ws.Range('B283').Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range('C283').Value = 'Arica y Parinacota'
$ws.Range('D283').Value = 44988
$ws.Range('E283').Value = 15
$ws.Range('F283').Value = 'Fruta'
$ws.Range('G283').Value = 100108
$ws.Range('H283').Value = 'Tropicales y subtropicales'
$ws.Range('I283').Value = 100108006
$ws.Range('J283').Value = 'Plátano'
$ws.Range('K283').Value = 'Sin especificar'
$ws.Range('L283').Value = 'Verde'
$ws.Range('M283').Value = 300
$ws.Range('N283').Value = 23000
$ws.Range('O283').Value = 24000
$ws.Range('P283').Value = 23500
$ws.Range('Q283').Value = '$/caja 20 kilos'
$ws.Range('R283').Value = 'Ecuador'
$ws.Range('S283').Value = 1175
$ws.Range('T283').Value = 20
